$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 2).Value = 1
    $ws.Cells.Item($row, 3).Value = "2025-04-04 13:22:34"
}
